$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 3
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()
# Row 15
$ws.Range("H15").Value = 4174.8647
$ws.Range("I15").Value = 4174.8647
$ws.Range("K15").Value = 12524.5941
$ws.Range("M15").Value = -12355.5941
# Row 98
$ws.Range("H98").Value = 4153.8335
$ws.Range("I98").Value = 4153.8335
$ws.Range("K98").Value = 4153.8335
$ws.Range("M98").Value = -2655.8335
# Row 102
$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()
# Row 112
$ws.Range("H112").Value = 1968
$ws.Range("I112").Value = 1229
$ws.Range("J112").Value = 2284.7144
$ws.Range("K112").Value = 3687
$ws.Range("L112").Value = 6854.1432
$ws.Range("M112").Value = -2579
$ws.Range("N112").Value = -9070.143199999999
# Row 122
$ws.Range("H122").Value = 4153.8335
$ws.Range("I122").Value = 4153.8335
$ws.Range("K122").Value = 12461.5005
$ws.Range("M122").Value = -10011.5005
# Row 132
$ws.Range("H132").Value = 5252.8
$ws.Range("I132").Value = 6358.864
$ws.Range("K132").Value = 19076.592
$ws.Range("M132").Value = -16546.592
# Row 137
$ws.Range("H137").Value = 4955
$ws.Range("I137").Value = 2409.1
$ws.Range("J137").Value = 8137.375
$ws.Range("K137").Value = 7227.299999999999
$ws.Range("L137").Value = 24412.125
$ws.Range("M137").Value = -4677.299999999999
$ws.Range("N137").Value = -29512.125
# Row 138
$ws.Range("H138").Value = 6733.1177
$ws.Range("J138").Value = 5917.0645
$ws.Range("L138").Value = 17751.1935
$ws.Range("N138").Value = -28031.1935
# Row 141
$ws.Range("H141").Value = 1879.5
$ws.Range("I141").Value = 1879.5
$ws.Range("K141").Value = 5638.5
$ws.Range("M141").Value = -458.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 4279.9653
$ws.Range("I61").Value = 3823.6365
$ws.Range("J61").Value = 5714.143
$ws.Range("K61").Value = 3823.6365
$ws.Range("L61").Value = 5714.143
$ws.Range("M61").Value = -3611.6365
$ws.Range("N61").Value = -6138.143
# Row 74
$ws.Range("H74").Value = 34677028
$ws.Range("I74").Value = 329289.12
$ws.Range("K74").Value = 329289.12
$ws.Range("M74").Value = -328415.12
# Row 77
$ws.Range("H77").Value = 34677028
$ws.Range("I77").Value = 329289.12
$ws.Range("K77").Value = 1646445.6
$ws.Range("M77").Value = -1642077.6
# Row 136
$ws.Range("H136").Value = 4279.9653
$ws.Range("I136").Value = 3823.6365
$ws.Range("J136").Value = 5714.143
$ws.Range("K136").Value = 11470.9095
$ws.Range("L136").Value = 17142.429
$ws.Range("M136").Value = -8920.9095
$ws.Range("N136").Value = -22242.429

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 15629551
$ws.Range("I20").Value = 19236084
$ws.Range("J20").Value = 1236.5
$ws.Range("K20").Value = 19236084
$ws.Range("L20").Value = 1236.5
$ws.Range("M20").Value = -19235837
$ws.Range("N20").Value = -1730.5
# Row 132
$ws.Range("H132").Value = 82093.375
$ws.Range("J132").Value = 82093.375
$ws.Range("L132").Value = 82093.375
$ws.Range("N132").Value = -92213.375
# Row 134
$ws.Range("H134").Value = 1623.1515
$ws.Range("I134").Value = 1288.1072
$ws.Range("K134").Value = 3864.3216
$ws.Range("M134").Value = -1329.3216

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 946.55554
$ws.Range("I22").Value = 859.8
$ws.Range("J22").Value = 1055
$ws.Range("K22").Value = 859.8
$ws.Range("L22").Value = 1055
$ws.Range("M22").Value = -509.8
$ws.Range("N22").Value = -1755
# Row 31
$ws.Range("H31").Value = 3157.1128
$ws.Range("I31").Value = 1245.5834
$ws.Range("J31").Value = 3615.88
$ws.Range("K31").Value = 1245.5834
$ws.Range("L31").Value = 3615.88
$ws.Range("M31").Value = -950.5834
$ws.Range("N31").Value = -4205.88
# Row 34
$ws.Range("H34").Value = 3157.1128
$ws.Range("I34").Value = 1245.5834
$ws.Range("J34").Value = 3615.88
$ws.Range("K34").Value = 1245.5834
$ws.Range("L34").Value = 3615.88
$ws.Range("M34").Value = -1043.5834
$ws.Range("N34").Value = -4019.88
# Row 41
$ws.Range("H41").Value = 55032.5
$ws.Range("J41").Value = 55065
$ws.Range("L41").Value = 55065
$ws.Range("N41").Value = -55921
# Row 45
$ws.Range("H45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("N45").ClearContents()
# Row 58
$ws.Range("H58").Value = 2523.8333
$ws.Range("I58").Value = 2026.7222
$ws.Range("K58").Value = 2026.7222
$ws.Range("M58").Value = -1823.7222
# Row 93
$ws.Range("H93").Value = 21562
$ws.Range("I93").Value = 21562
$ws.Range("K93").Value = 21562
$ws.Range("M93").Value = -19690
# Row 136
$ws.Range("H136").Value = 2523.8333
$ws.Range("I136").Value = 2026.7222
$ws.Range("K136").Value = 6080.1666
$ws.Range("M136").Value = -3530.1666

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 7
$ws.Range("H7").Value = 30304512
$ws.Range("I7").Value = 33334918
$ws.Range("J7").Value = 450
$ws.Range("K7").Value = 100004754
$ws.Range("L7").Value = 1350
$ws.Range("M7").Value = -100004642
$ws.Range("N7").Value = -1574
# Row 131
$ws.Range("H131").Value = 7807.4062
$ws.Range("J131").Value = 2103.682
$ws.Range("L131").Value = 6311.045999999999
$ws.Range("N131").Value = -16391.046
# Row 137
$ws.Range("H137").Value = 50000
$ws.Range("I137").Value = 50000
$ws.Range("K137").Value = 150000
$ws.Range("M137").Value = -144900

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Range("H122").Value = 21279510
$ws.Range("I122").Value = 2546.9333
$ws.Range("J122").Value = 58827092
$ws.Range("K122").Value = 7640.7999
$ws.Range("L122").Value = 176481276
$ws.Range("M122").Value = -5190.7999
$ws.Range("N122").Value = -176486176
# Row 132
$ws.Range("H132").Value = 1955.762
$ws.Range("I132").Value = 1494.1305
$ws.Range("K132").Value = 4482.3915
$ws.Range("M132").Value = -1952.3915

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 48383.332
$ws.Range("I40").Value = 88447.71000000001
$ws.Range("K40").Value = 88447.71000000001
$ws.Range("M40").Value = -88311.71000000001
# Row 96
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()
# Row 136
$ws.Range("H136").Value = 6649.273
$ws.Range("I136").Value = 6402.6665
$ws.Range("K136").Value = 19207.9995
$ws.Range("M136").Value = -16657.9995

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 1599.9348
$ws.Range("I132").Value = 1368.0256
$ws.Range("J132").Value = 2892
$ws.Range("K132").Value = 2892
$ws.Range("L132").Value = 8676
$ws.Range("M132").Value = -1574.0768
$ws.Range("N132").Value = -13736
